$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 22.481209
$ws.Range("H2").Value = 67.443627
$ws.Range("I2").Value = 0.1656226259370683
$ws.Range("J2").Value = 0.166106832923046
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.787414
$ws.Range("N2").Value = 11.362242
$ws.Range("O2").Value = 0.5877125485801681
$ws.Range("P2").Value = 0.587712548580168
$ws.Range("Q2").Value = 85.14564570352601
$ws.Range("R2").Value = 766.310811331734
$ws.Range("S2").Value = 0.09733849559201428
$ws.Range("T2").Value = 0.0976230701137835

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 22.481209
$ws.Range("H3").Value = 67.443627
$ws.Range("I3").Value = 0.1656226259370683
$ws.Range("J3").Value = 0.166106832923046
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.656916666666667
$ws.Range("N3").Value = 7.97075
$ws.Range("O3").Value = 0.412287451419832
$ws.Range("P3").Value = 0.4122874514198319
$ws.Range("Q3").Value = 59.73069887891668
$ws.Range("R3").Value = 537.57628991025
$ws.Range("S3").Value = 0.06828413034505407
$ws.Range("T3").Value = 0.06848376280926245

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 29.50180766666667
$ws.Range("H4").Value = 88.50542300000001
$ws.Range("I4").Value = 0.2173444878184117
$ws.Range("J4").Value = 0.2179799065528387
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.787414
$ws.Range("N4").Value = 11.362242
$ws.Range("O4").Value = 0.5877125485801681
$ws.Range("P4").Value = 0.587712548580168
$ws.Range("Q4").Value = 111.7355593820407
$ws.Range("R4").Value = 1005.620034438366
$ws.Range("S4").Value = 0.12773608285561
$ws.Range("T4").Value = 0.1281095264194357

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.50180766666667
$ws.Range("H5").Value = 88.50542300000001
$ws.Range("I5").Value = 0.2173444878184117
$ws.Range("J5").Value = 0.2179799065528387
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.656916666666667
$ws.Range("N5").Value = 7.97075
$ws.Range("O5").Value = 0.412287451419832
$ws.Range("P5").Value = 0.4122874514198319
$ws.Range("Q5").Value = 78.38384448636111
$ws.Range("R5").Value = 705.45460037725
$ws.Range("S5").Value = 0.08960840496280167
$ws.Range("T5").Value = 0.089870380133403

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 42.765269
$ws.Range("H6").Value = 128.295807
$ws.Range("I6").Value = 0.3150585073376215
$ws.Range("J6").Value = 0.3159795984589671
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.787414
$ws.Range("N6").Value = 11.362242
$ws.Range("O6").Value = 0.5877125485801681
$ws.Range("P6").Value = 0.587712548580168
$ws.Range("Q6").Value = 161.969778524366
$ws.Range("R6").Value = 1457.728006719294
$ws.Range("S6").Value = 0.1851638382992571
$ws.Range("T6").Value = 0.1857051751096576

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 42.765269
$ws.Range("H7").Value = 128.295807
$ws.Range("I7").Value = 0.3150585073376215
$ws.Range("J7").Value = 0.3159795984589671
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.656916666666667
$ws.Range("N7").Value = 7.97075
$ws.Range("O7").Value = 0.412287451419832
$ws.Range("P7").Value = 0.4122874514198319
$ws.Range("Q7").Value = 113.6237559605833
$ws.Range("R7").Value = 1022.61380364525
$ws.Range("S7").Value = 0.1298946690383644
$ws.Range("T7").Value = 0.1302744233493094

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 39.80222300000001
$ws.Range("H8").Value = 119.406669
$ws.Range("I8").Value = 0.2932292783449848
$ws.Range("J8").Value = 0.2940865504976542
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.787414
$ws.Range("N8").Value = 11.362242
$ws.Range("O8").Value = 0.5877125485801681
$ws.Range("P8").Value = 0.587712548580168
$ws.Range("Q8").Value = 150.747496621322
$ws.Range("R8").Value = 1356.727469591898
$ws.Range("S8").Value = 0.1723345264944545
$ws.Range("T8").Value = 0.1728383560961266

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 39.80222300000001
$ws.Range("H9").Value = 119.406669
$ws.Range("I9").Value = 0.2932292783449848
$ws.Range("J9").Value = 0.2940865504976542
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.656916666666667
$ws.Range("N9").Value = 7.97075
$ws.Range("O9").Value = 0.412287451419832
$ws.Range("P9").Value = 0.4122874514198319
$ws.Range("Q9").Value = 105.7511896590833
$ws.Range("R9").Value = 951.7607069317501
$ws.Range("S9").Value = 0.1208947518505303
$ws.Range("T9").Value = 0.1212481944015276

$ws.Range("E10").Value = 2
$ws.Range("G10").Value = 1.1870385
$ws.Range("H10").Value = 2.374077
$ws.Range("I10").Value = 0.00874510056191367
$ws.Range("J10").Value = 0.005847111567493934
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.787414
$ws.Range("N10").Value = 11.362242
$ws.Range("O10").Value = 0.5877125485801681
$ws.Range("P10").Value = 0.587712548580168
$ws.Range("Q10").Value = 4.495806233439001
$ws.Range("R10").Value = 26.974837400634
$ws.Range("S10").Value = 0.005139605338832143
$ws.Range("T10").Value = 0.003436420841164441

$ws.Range("E11").Value = 2
$ws.Range("G11").Value = 1.1870385
$ws.Range("H11").Value = 2.374077
$ws.Range("I11").Value = 0.00874510056191367
$ws.Range("J11").Value = 0.005847111567493934
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.656916666666667
$ws.Range("N11").Value = 7.97075
$ws.Range("O11").Value = 0.412287451419832
$ws.Range("P11").Value = 0.4122874514198319
$ws.Range("Q11").Value = 3.153862374625
$ws.Range("R11").Value = 18.92317424775
$ws.Range("S11").Value = 0.003605495223081528
$ws.Range("T11").Value = 0.002410690726329493

